# Update Lab/Hwk schedule links on the "overview" sheet.
# The Lab (col G) and Hwk (col H) links are being shifted down to new rows,
# and a new "Hwk 5" link is added at H12.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("overview")

# --- Row 3 ---
# Hwk 1 link removed from H3 (it moves to H4)
$ws.Range("H3").ClearContents()

# --- Row 4 ---
# Lab 1 link removed from G4 (it moves to G5)
$ws.Range("G4").ClearContents()
# Hwk 1 link now placed in H4 (was previously Hwk 2, which moves to H6)
$ws.Range("H4").Value = "[Hwk 1](https://tinyurl.com/demog180-fa23-hwk01)"

# --- Row 5 ---
# Lab 1 link now placed in G5 (was previously Lab 2, which moves to G7)
$ws.Range("G5").Value = "[Lab 1](https://datahub.berkeley.edu/hub/user-redirect/git-pull?repo=https%3A%2F%2Fgithub.com%2Fdfeehan%2Fdemog180-fa2023&urlpath=tree%2Fdemog180-fa2023%2Flab%2Flab1%2Flab_complete_network_data.ipynb&branch=main)"
# Hwk 3 link removed from H5 (it moves to H8)
$ws.Range("H5").ClearContents()

# --- Row 6 ---
# Lab 3 link removed from G6 (it moves to G9)
$ws.Range("G6").ClearContents()
# Hwk 2 link now placed in H6 (was previously Hwk 4, which moves to H10)
$ws.Range("H6").Value = "[Hwk 2](https://datahub.berkeley.edu/hub/user-redirect/git-pull?repo=https%3A%2F%2Fgithub.com%2Fdfeehan%2Fdemog180-fa2023&urlpath=tree%2Fdemog180-fa2023%2Fhwk%2Fhwk02%2Fhwk_clusteringcoef.ipynb&branch=main)"

# --- Row 7 ---
# Lab 2 link now placed in G7 (was previously Lab 4, which moves to G11)
$ws.Range("G7").Value = "[Lab 2](https://datahub.berkeley.edu/hub/user-redirect/git-pull?repo=https%3A%2F%2Fgithub.com%2Fdfeehan%2Fdemog180-fa2023&urlpath=tree%2Fdemog180-fa2023%2Flab%2Flab2%2Flab2_personal_networks.ipynb&branch=main)"

# --- Row 8 ---
# Hwk 3 link added at H8
$ws.Range("H8").Value = "[Hwk 3](https://datahub.berkeley.edu/hub/user-redirect/git-pull?repo=https%3A%2F%2Fgithub.com%2Fdfeehan%2Fdemog180-fa2023&urlpath=tree%2Fdemog180-fa2023%2Fhwk%2Fhwk03%2Fhwk03_personal_networks.ipynb&branch=main)"

# --- Row 9 ---
# Lab 3 link added at G9
$ws.Range("G9").Value = "[Lab 3](https://datahub.berkeley.edu/hub/user-redirect/git-pull?repo=https%3A%2F%2Fgithub.com%2Fdfeehan%2Fdemog180-fa2023&urlpath=tree%2Fdemog180-fa2023%2Flab%2Flab3%2Flab03_homophily.ipynb&branch=main)"

# --- Row 10 ---
# Hwk 4 link added at H10
$ws.Range("H10").Value = "[Hwk 4](https://datahub.berkeley.edu/hub/user-redirect/git-pull?repo=https%3A%2F%2Fgithub.com%2Fdfeehan%2Fdemog180-fa2023&urlpath=tree%2Fdemog180-fa2023%2Fhwk%2Fhwk04%2Fhw04_balance_smallworlds.ipynb&branch=main)"

# --- Row 11 ---
# Lab 4 link added at G11
$ws.Range("G11").Value = "[Lab 4](https://datahub.berkeley.edu/hub/user-redirect/git-pull?repo=https%3A%2F%2Fgithub.com%2Fdfeehan%2Fdemog180-fa2023&urlpath=tree%2Fdemog180-fa2023%2Flab%2Flab%25204%2Flab04_affiliation_networks.ipynb&branch=main)"

# --- Row 12 ---
# New: Hwk 5 link added at H12
$ws.Range("H12").Value = "[Hwk 5](https://datahub.berkeley.edu/hub/user-redirect/git-pull?repo=https%3A%2F%2Fgithub.com%2Fdfeehan%2Fdemog180-fa2023&urlpath=tree%2Fdemog180-fa2023%2Fhwk%2Fhwk05%2Fhwk05_small_worlds.ipynb&branch=main)"
